$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.453.31'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').Value = '3.313.37'
$ws.Range('E3').Value = '  +1.31%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '186.11'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '576.92'
$ws.Range('E6').Value = '  -1.17%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  -0.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  +1.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.409'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '3.888.05'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('E13').Value = '  -0.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.46'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('D15').Value = '67.654.90'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = '3.312.58'
$ws.Range('E17').Value = '  -0.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '444.13'
$ws.Range('E18').Value = '  +6.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.69'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.76'
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.04'
$ws.Range('E22').Value = '  +3.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.517'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('D25').Value = '3.455.94'
$ws.Range('E25').Value = '  +1.24%  '
$ws.Range('E26').Value = '  +0.91%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.04'
$ws.Range('E28').Value = '  -4.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.94'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.34'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.25'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  +4.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.98'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.28'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.48'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('D42').Value = '2.760.02'
$ws.Range('E42').Value = '  +3.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.28'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.85'
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.21'
$ws.Range('E45').Value = '  -1.58%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0672'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '326.45'
$ws.Range('E48').Value = '  -3.23%  '
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.992'
$ws.Range('E50').Value = '  +1.88%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.27'
$ws.Range('E51').Value = '  +1.66%  '
